# Updates the cryptos price/volume table (GitHub Actions-style scrape refresh).
# Price (col D) and Volume(1h) (col E) cells are stored as plain text in the
# sheet, not numbers. For D-column values that look numeric to Excel's
# auto-detection (e.g. "250.18"), we force the cell to Text format before
# assigning the string, then ClearFormats() to drop the temporary "@" number
# format again (restoring the cell to its original unstyled state) while the
# text datatype sticks. Values that already contain non-numeric punctuation
# (two dots, subscript digits, etc.) are left untouched since Excel keeps
# those as text automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.776.50'
$ws.Range('E2').Value = '  +4.11%  '
$ws.Range('D3').Value = '1.914.15'
$ws.Range('E3').Value = '  +1.56%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.18'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.703'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.88%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '46.49'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +7.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.373'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +5.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.52'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +9.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0763'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.67%  '
$ws.Range('E12').Value = '  +1.80%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.65'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +8.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.814'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +5.62%  '
$ws.Range('D15').Value = '2.194.95'
$ws.Range('E15').Value = '  +1.70%  '
$ws.Range('E16').Value = '  +3.85%  '
$ws.Range('D17').Value = '1.914.69'
$ws.Range('E17').Value = '  +1.23%  '
$ws.Range('D18').Value = '36.758.02'
$ws.Range('E18').Value = '  +4.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '75.02'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.39%  '
$ws.Range('D20').Value = '0.0₃0859'
$ws.Range('E20').Value = '  +3.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '252.13'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.39'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +4.43%  '
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('E24').Value = '  +0.50%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.21'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.87'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.81'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.95%  '
$ws.Range('E29').Value = '  +2.58%  '
$ws.Range('E30').Value = '  +1.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.59'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +6.56%  '
$ws.Range('E32').Value = '  +4.13%  '
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0911'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +24.31%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.34'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +3.22%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.90'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.44%  '
$ws.Range('B36').Value = 'BinanceUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.53'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +5.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.38'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +54.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.876'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.68%  '
$ws.Range('E40').Value = '  +2.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '105.39'
$ws.Range('D41').ClearFormats()
$ws.Range('E42').Value = '  +4.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.74'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.86'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +19.58%  '
$ws.Range('E45').Value = '  +2.19%  '
$ws.Range('D46').Value = '1.350.00'
$ws.Range('E46').Value = '  +2.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.39'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0814'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.33%  '
$ws.Range('E49').Value = '  +2.35%  '
$ws.Range('E50').Value = '  +2.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.37'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.99%  '
